$p = $ppt.ActivePresentation
$s = $p.Slides.Item(13)

# Locate the content placeholder shape (id=5, "Inhaltsplatzhalter 2") robustly
# by its shape Id rather than a fragile positional index.
$sp = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $cand = $s.Shapes.Item($i)
    if ($cand.Id -eq 5) {
        $sp = $cand
        break
    }
}
if ($sp -eq $null) {
    $sp = $s.Shapes.Item(2)
}
$tr = $sp.TextFrame.TextRange

# --- Open Point 1: bold "that anydata modeled nodes can be validated " ---
$full = $tr.Text
$span1 = "that anydata modeled nodes can be validated "
$idx1 = $full.IndexOf($span1)
if ($idx1 -ge 0) {
    $sub1 = $tr.Characters($idx1 + 1, $span1.Length)
    $sub1.Font.Bold = 1
}

# --- Open Point 2: bold "validates the content of anydata nodes " ---
$full = $tr.Text
$span2 = "validates the content of anydata nodes "
$idx2 = $full.IndexOf($span2)
if ($idx2 -ge 0) {
    $sub2 = $tr.Characters($idx2 + 1, $span2.Length)
    $sub2.Font.Bold = 1
}

# --- Open Point 3: "JSON are missing" -> "JSON and CBOR are missing" ---
# Target the whole run's text (matching its exact boundaries) so PowerPoint
# updates it in place instead of splitting into extra runs.
$full = $tr.Text
$span3 = ".  However, specifications for encoding in JSON are missing "
$idx3 = $full.IndexOf($span3)
if ($idx3 -ge 0) {
    $sub3 = $tr.Characters($idx3 + 1, $span3.Length)
    $sub3.Text = ".  However, specifications for encoding in JSON and CBOR are missing "
}
